$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix the body of the hello() rule: it now delegates to hello3()
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = "return hello3() ;"

# ---------------------------------------------------------------------------
# 2. New rule table "Rules String hello3()" in B8:D13
#    Row 8 is the table header (same look as the B3/B5 header bars).
# ---------------------------------------------------------------------------
$header = $ws.Range("B8:D8")
$header.Merge()
$header.Value = "Rules String hello3()"
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108

# Rows 9-13: each one is its own merged B:D row, with a thin box border
# (left edge only on column B, right edge only on column D, top+bottom thin
# on every row so each row is visually separated) and centered text.
function Set-BoxRow($rowNum, $text) {
    $rng = $ws.Range("B$rowNum" + ":D$rowNum")
    $rng.Merge()
    $rng.Value = $text
    $rng.HorizontalAlignment = -4108
    $rng.Borders.LineStyle = 1
    $rng.Borders.Item(11).LineStyle = 0
    $rng.Borders.Item(12).LineStyle = 0
}

Set-BoxRow 9  "RET1"
Set-BoxRow 10 "res"
Set-BoxRow 11 "String res"
Set-BoxRow 12 "From dependency"
# Row 13 holds literal text that starts with "=" -- the leading apostrophe
# forces it to stay text instead of being parsed as a formula (quote prefix).
Set-BoxRow 13 "'=return helloFromDependency() ;"

# ---------------------------------------------------------------------------
# 3. "Environment" / dependency declaration block in B15:C16
# ---------------------------------------------------------------------------
$envHeader = $ws.Range("B15:C15")
$envHeader.Merge()
$envHeader.Value = "Environment"
$envHeader.HorizontalAlignment = -4108

$ws.Range("B16").Value = "dependency"
$ws.Range("C16").Value = "Module3_2"

# ---------------------------------------------------------------------------
# 4. View state: scroll so row 4 is at the top and select B14
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("B14").Select()
